$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before D (shifts D:K -> E:L), matching the workbook-wide
# restructuring that added a newest reporting period as the first data column.
$ws.Columns("D:D").Insert()

# Carry over number formatting/styles from the old first data column (now column E)
# onto the freshly inserted column D so dates/numbers render the same way.
$ws.Range("E7:E102").Copy()
$ws.Range("D7:D102").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Columns("D:D").ColumnWidth = $ws.Columns("E:E").ColumnWidth

# Populate new column D with the new period values
$newValues = @{
    7 = 43465
    8 = 213100
    13 = 0
    14 = 0
    15 = -600
    17 = 44000
    18 = 169100
    20 = -88900
    21 = 86400
    22 = 0
    23 = 80300
    24 = 16900
    25 = 0
    26 = 63400
    27 = 62800
    28 = 0
    29 = 0
    30 = 0
    31 = 0
    32 = 88900
    33 = 62800
    34 = 0
    35 = 62800
    38 = 43465
    41 = 205200
    42 = 16700
    43 = 0
    44 = 0
    45 = 0
    46 = 0
    48 = 49200
    49 = 138200
    50 = 0
    51 = 0
    52 = 17200
    53 = 0
    54 = 5806100
    57 = 0
    58 = 0
    59 = 0
    60 = 0
    61 = 105000
    62 = 0
    63 = 0
    64 = 0
    65 = 0
    66 = 5182400
    68 = 0
    69 = 0
    70 = 0
    71 = 0
    72 = 116900
    73 = 0
    74 = 0
    75 = 0
    76 = 623700
    77 = 0
    80 = 43465
    81 = 62800
    83 = 6100
    84 = 0
    85 = 0
    86 = 0
    87 = 0
    88 = 0
    89 = 79400
    91 = -5500
    92 = 0
    93 = 0
    94 = -342600
    96 = -21300
    97 = 0
    98 = 0
    99 = 0
    100 = 328800
    101 = 0
    102 = 65700
}
foreach ($row in $newValues.Keys) {
    $ws.Range("D" + $row).Value2 = $newValues[$row]
}

# These rows hold "NA" (text) markers in the new column D as well
$naRows = @(9, 10, 12)
foreach ($row in $naRows) {
    $ws.Range("D" + $row).Value2 = "NA"
}

# Row 47 (Deferred Long Term Asset Charges) changed shape: new D gets a value,
# E:J flip from 0 to "NA", and K (shifted from old J) stays 0.
$ws.Range("D47").Value2 = 15900
$ws.Range("E47:J47").Value2 = "NA"

Write-Host "Applied LBAI financials update"
